# Weekly crime-data refresh: advance the report one week
# (Volume/Number bump, week-covering dates, and the Week-to-Date / 28-Day /
#  Year-to-Date / 2-Year crime-count table for rows 15-27).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: "Volume 30   Number  37" -> "...  38", and the week span ---
$ws.Range("A8").Value = "Volume 30   Number  38"
$ws.Range("C9").Value = "Report Covering the Week  9/18/2023  Through  9/24/2023"

# Helper: a couple of cells (D15/E15/C16/D17/E17/D26/E26) flip from a plain
# number to the sheet's "no data" placeholder text ("0" or "***.*"). Those
# placeholder cells elsewhere on the sheet (e.g. C15) use cell style 14
# (General number format, right aligned) instead of the numeric style the
# cell currently has, so after poking in the text we re-stamp the format
# from a known-good style-14 cell (C15) via a formats-only paste.
function Set-PlaceholderText {
    param($cellRef, $text)
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range("C15").Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null
}

# Row 15 - Rape
Set-PlaceholderText "D15" "0"
Set-PlaceholderText "E15" "***.*"

# Row 16 - Robbery
Set-PlaceholderText "C16" "0"
$ws.Range("E16").Value = -100
$ws.Range("J16").Value = 16
$ws.Range("K16").Value = -6.25
$ws.Range("L16").Value = 150
$ws.Range("M16").Value = -21.052631578947
$ws.Range("N16").Value = -65.116279069767

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 2
Set-PlaceholderText "D17" "0"
Set-PlaceholderText "E17" "***.*"
$ws.Range("F17").Value = 5
$ws.Range("G17").Value = 3
$ws.Range("H17").Value = 66.666666666666
$ws.Range("I17").Value = 68
$ws.Range("K17").Value = 112.5
$ws.Range("L17").Value = 151.851851851852
$ws.Range("M17").Value = 106.060606060606
$ws.Range("N17").Value = -17.073170731707

# Row 18 - Burglary
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -62.5
$ws.Range("I18").Value = 42
$ws.Range("J18").Value = 27
$ws.Range("K18").Value = 55.555555555555
$ws.Range("L18").Value = 121.052631578947
$ws.Range("M18").Value = -46.153846153846
$ws.Range("N18").Value = -81.415929203539

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 0
$ws.Range("I19").Value = 207
$ws.Range("J19").Value = 205
$ws.Range("K19").Value = 0.975609756097
$ws.Range("L19").Value = 71.074380165289
$ws.Range("M19").Value = 102.941176470588
$ws.Range("N19").Value = 23.952095808383

# Row 20 - G.L.A.
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = -75
$ws.Range("I20").Value = 61
$ws.Range("J20").Value = 85
$ws.Range("K20").Value = -28.235294117647
$ws.Range("L20").Value = 79.411764705882
$ws.Range("M20").Value = 125.925925925926
$ws.Range("N20").Value = -88.403041825095

# Row 21 - TOTAL
$ws.Range("D21").Value = 12
$ws.Range("E21").Value = -8.333333333333
$ws.Range("F21").Value = 26
$ws.Range("G21").Value = 55
$ws.Range("H21").Value = -52.727272727272
$ws.Range("I21").Value = 395
$ws.Range("J21").Value = 367
$ws.Range("K21").Value = 7.629427792915
$ws.Range("L21").Value = 88.995215311004
$ws.Range("M21").Value = 51.340996168582
$ws.Range("N21").Value = -62.309160305343

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = -42.857142857142
$ws.Range("G24").Value = 44
$ws.Range("H24").Value = -22.727272727272
$ws.Range("I24").Value = 345
$ws.Range("J24").Value = 334
$ws.Range("K24").Value = 3.293413173652
$ws.Range("L24").Value = 79.6875
$ws.Range("M24").Value = -13.316582914572

# Row 25 - Misd. Assault
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 66.666666666666
$ws.Range("F25").Value = 12
$ws.Range("G25").Value = 16
$ws.Range("H25").Value = -25
$ws.Range("I25").Value = 149
$ws.Range("J25").Value = 123
$ws.Range("K25").Value = 21.138211382113
$ws.Range("L25").Value = 47.524752475247
$ws.Range("M25").Value = -9.696969696969

# Row 26 - UCR Rape*
Set-PlaceholderText "D26" "0"
Set-PlaceholderText "E26" "***.*"

# Row 27 - Other Sex Crimes
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 14
$ws.Range("K27").Value = -6.666666666666
$ws.Range("L27").Value = 55.555555555555
